$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2247043333333333
$ws.Range("H2").Value = 0.674113
$ws.Range("I2").Value = 0.2389319335355998
$ws.Range("J2").Value = 0.2389319335355999
$ws.Range("O2").Value = 0.06522509891308133
$ws.Range("P2").Value = 0.06522509891308133
$ws.Range("Q2").Value = 0.04602289293311111
$ws.Range("R2").Value = 0.414206036398
$ws.Range("S2").Value = 0.01558435899835327
$ws.Range("T2").Value = 0.01558435899835328
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2247043333333333
$ws.Range("H3").Value = 0.674113
$ws.Range("I3").Value = 0.2389319335355998
$ws.Range("J3").Value = 0.2389319335355999
$ws.Range("M3").Value = 0.6481333333333333
$ws.Range("N3").Value = 1.9444
$ws.Range("O3").Value = 0.2064033004146749
$ws.Range("P3").Value = 0.2064033004146749
$ws.Range("Q3").Value = 0.1456383685777778
$ws.Range("R3").Value = 1.3107453172
$ws.Range("S3").Value = 0.04931633965620755
$ws.Range("T3").Value = 0.04931633965620755
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2247043333333333
$ws.Range("H4").Value = 0.674113
$ws.Range("I4").Value = 0.2389319335355998
$ws.Range("J4").Value = 0.2389319335355999
$ws.Range("M4").Value = 1.888205
$ws.Range("N4").Value = 5.664615
$ws.Range("O4").Value = 0.6013141491351952
$ws.Range("P4").Value = 0.6013141491351952
$ws.Range("Q4").Value = 0.4242878457216666
$ws.Range("R4").Value = 3.818590611494999
$ws.Range("S4").Value = 0.1436731523151862
$ws.Range("T4").Value = 0.1436731523151863
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2247043333333333
$ws.Range("H5").Value = 0.674113
$ws.Range("I5").Value = 0.2389319335355998
$ws.Range("J5").Value = 0.2389319335355999
$ws.Range("M5").Value = 0.398977
$ws.Range("N5").Value = 1.196931
$ws.Range("O5").Value = 0.1270574515370486
$ws.Range("P5").Value = 0.1270574515370486
$ws.Range("Q5").Value = 0.08965186080033331
$ws.Range("R5").Value = 0.8068667472029999
$ws.Range("S5").Value = 0.03035808256585278
$ws.Range("T5").Value = 0.03035808256585279
$ws.Range("G6").Value = 0.4451493333333333
$ws.Range("H6").Value = 1.335448
$ws.Range("I6").Value = 0.4733348456063742
$ws.Range("J6").Value = 0.4733348456063743
$ws.Range("O6").Value = 0.06522509891308133
$ws.Range("P6").Value = 0.06522509891308133
$ws.Range("Q6").Value = 0.0911734090897778
$ws.Range("R6").Value = 0.820560681808
$ws.Range("S6").Value = 0.03087331212368384
$ws.Range("T6").Value = 0.03087331212368384
$ws.Range("G7").Value = 0.4451493333333333
$ws.Range("H7").Value = 1.335448
$ws.Range("I7").Value = 0.4733348456063742
$ws.Range("J7").Value = 0.4733348456063743
$ws.Range("M7").Value = 0.6481333333333333
$ws.Range("N7").Value = 1.9444
$ws.Range("O7").Value = 0.2064033004146749
$ws.Range("P7").Value = 0.2064033004146749
$ws.Range("Q7").Value = 0.2885161212444445
$ws.Range("R7").Value = 2.5966450912
$ws.Range("S7").Value = 0.09769787433442623
$ws.Range("T7").Value = 0.09769787433442624
$ws.Range("G8").Value = 0.4451493333333333
$ws.Range("H8").Value = 1.335448
$ws.Range("I8").Value = 0.4733348456063742
$ws.Range("J8").Value = 0.4733348456063743
$ws.Range("M8").Value = 1.888205
$ws.Range("N8").Value = 5.664615
$ws.Range("O8").Value = 0.6013141491351952
$ws.Range("P8").Value = 0.6013141491351952
$ws.Range("Q8").Value = 0.8405331969466666
$ws.Range("R8").Value = 7.564798772519999
$ws.Range("S8").Value = 0.2846229399418359
$ws.Range("T8").Value = 0.284622939941836
$ws.Range("G9").Value = 0.4451493333333333
$ws.Range("H9").Value = 1.335448
$ws.Range("I9").Value = 0.4733348456063742
$ws.Range("J9").Value = 0.4733348456063743
$ws.Range("M9").Value = 0.398977
$ws.Range("N9").Value = 1.196931
$ws.Range("O9").Value = 0.1270574515370486
$ws.Range("P9").Value = 0.1270574515370486
$ws.Range("Q9").Value = 0.1776043455653333
$ws.Range("R9").Value = 1.598439110088
$ws.Range("S9").Value = 0.06014071920642826
$ws.Range("T9").Value = 0.06014071920642827
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.2705996666666666
$ws.Range("H10").Value = 0.8117989999999999
$ws.Range("I10").Value = 0.2877332208580259
$ws.Range("J10").Value = 0.2877332208580259
$ws.Range("O10").Value = 0.06522509891308133
$ws.Range("P10").Value = 0.06522509891308133
$ws.Range("Q10").Value = 0.05542296092822222
$ws.Range("R10").Value = 0.498806648354
$ws.Range("S10").Value = 0.01876742779104421
$ws.Range("T10").Value = 0.01876742779104422
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.2705996666666666
$ws.Range("H11").Value = 0.8117989999999999
$ws.Range("I11").Value = 0.2877332208580259
$ws.Range("J11").Value = 0.2877332208580259
$ws.Range("M11").Value = 0.6481333333333333
$ws.Range("N11").Value = 1.9444
$ws.Range("O11").Value = 0.2064033004146749
$ws.Range("P11").Value = 0.2064033004146749
$ws.Range("Q11").Value = 0.1753846639555555
$ws.Range("R11").Value = 1.5784619756
$ws.Range("S11").Value = 0.05938908642404112
$ws.Range("T11").Value = 0.05938908642404113
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.2705996666666666
$ws.Range("H12").Value = 0.8117989999999999
$ws.Range("I12").Value = 0.2877332208580259
$ws.Range("J12").Value = 0.2877332208580259
$ws.Range("M12").Value = 1.888205
$ws.Range("N12").Value = 5.664615
$ws.Range("O12").Value = 0.6013141491351952
$ws.Range("P12").Value = 0.6013141491351952
$ws.Range("Q12").Value = 0.5109476435983332
$ws.Range("R12").Value = 4.598528792384999
$ws.Range("S12").Value = 0.173018056878173
$ws.Range("T12").Value = 0.1730180568781731
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.2705996666666666
$ws.Range("H13").Value = 0.8117989999999999
$ws.Range("I13").Value = 0.2877332208580259
$ws.Range("J13").Value = 0.2877332208580259
$ws.Range("M13").Value = 0.398977
$ws.Range("N13").Value = 1.196931
$ws.Range("O13").Value = 0.1270574515370486
$ws.Range("P13").Value = 0.1270574515370486
$ws.Range("Q13").Value = 0.1079630432076666
$ws.Range("R13").Value = 0.971667388869
$ws.Range("S13").Value = 0.03655864976476751
$ws.Range("T13").Value = 0.03655864976476752
